$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data Manipulation / Tableau cell changes from High to Medium
$ws.Range("C4").Value = "Medium"

# Caption below the table
$ws.Range("A13").Value = "Table 1. Technology Capability Ranked High to Low"
$ws.Range("A13").Font.Bold = $true

# Rename "Advanced Techniques" row label to "Advanced Modeling"
$ws.Range("A8").Value = "Advanced Modeling"

# Bold the header row (software names) and the row-label column (capability names)
$ws.Range("B2:D2").Font.Bold = $true
$ws.Range("A3:A11").Font.Bold = $true

# Update selection to match the saved view
$ws.Range("C9").Select()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
